{"js": "// Update the title: \"Results for attempt 1\" -> \"Results for attempt 17\"\nconst titleResults = context.document.body.search(\"Results for attempt 1\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\"Results for attempt 17\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Grab the two tables in the document: Quantitative (index 0) and Qualitative (index 1)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Quantitative table: Repetitions 148 -> 0, Duration 25.206612 -> 1.664161\nconst quantTable = tables.items[0];\nconst repetitionsCell = quantTable.getCell(1, 0);\nrepetitionsCell.body.insertText(\"0\", Word.InsertLocation.replace);\nconst durationCell = quantTable.getCell(1, 1);\ndurationCell.body.insertText(\"1.664161\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Qualitative table: Accuracy 148.0 -> 88.483343\nconst qualTable = tables.items[1];\nconst accuracyCell = qualTable.getCell(1, 0);\naccuracyCell.body.insertText(\"88.483343\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Title: \"Results for attempt 1\" -> \"Results for attempt 17\"\n$find = $d.Content.Find\n$find.Text = \"Results for attempt 1\"\n$find.Replacement.Text = \"Results for attempt 17\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Quantitative table (table 1): Repetitions 148 -> 0, Duration 25.206612 -> 1.664161\n$quantTable = $d.Tables.Item(1)\n$quantTable.Cell(2, 1).Range.Text = \"0\"\n$quantTable.Cell(2, 2).Range.Text = \"1.664161\"\n\n# 3. Qualitative table (table 2): Accuracy 148.0 -> 88.483343\n$qualTable = $d.Tables.Item(2)\n$qualTable.Cell(2, 1).Range.Text = \"88.483343\"\n"}
